# Fix for getWFTaskTable loop bug: append the newly-documented trial-table
# variables (bumpDir, bumpPhase, bumpTime, tgtOnTime) to the label sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "variable name" column (A) values first, then "description" column (B)
# values, so new shared-string entries are appended in the same order as
# the canonical workbook (all names, then all descriptions).
$ws.Range("A15").Value2 = "bumpDir"
$ws.Range("A16").Value2 = "bumpPhase"
$ws.Range("A17").Value2 = "bumpTime"
$ws.Range("A18").Value2 = "tgtOnTime"

$ws.Range("B15").Value2 = "direction of the bump in the room coordinate system"
$ws.Range("B16").Value2 = "what part of the trial was the bump in, e.g hold period bumps = 'H'"
$ws.Range("B17").Value2 = "time after trial start when bump was initiated"
$ws.Range("B18").Value2 = "time after trial start at which reach targets is presented"

# Match the author's final selection (one row below the last used row).
$ws.Range("B19").Select()

$wb.Save()
